$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 71

# Date and Week columns would be auto-converted to numbers/dates by Excel's
# smart input parsing (e.g. "2024-01-18" -> date serial, "02" -> 2), so mark
# them as Text before assigning so the literal string is preserved.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2024-01-18"

$ws.Range("B$row").Value = "18:18:54"
$ws.Range("C$row").Value = "Thursday"

$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "02"

$ws.Range("E$row").Value = 139068
$ws.Range("F$row").Value = 140168
$ws.Range("G$row").Value = 171513
$ws.Range("H$row").Value = 148765
$ws.Range("I$row").Value = -1
$ws.Range("J$row").Value = 121193
$ws.Range("K$row").Value = 223364
$ws.Range("L$row").Value = 254746
$ws.Range("M$row").Value = 185087
$ws.Range("N$row").Value = 110337
$ws.Range("O$row").Value = 41348
$ws.Range("P$row").Value = 30916
$ws.Range("Q$row").Value = 73545
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 42688
$ws.Range("T$row").Value = -1
